$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Acelga at Feria Lagunitas
# de Puerto Montt. Insert it as a new row 129 (pushing the existing rows
# 129-140 down to 130-141, same as the source diff).
$ws.Rows(129).Insert()

$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44578
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = 100112009
$ws.Range("G129").Value = "Acelga"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 40
$ws.Range("K129").Value = 10000
$ws.Range("L129").Value = 10000
$ws.Range("M129").Value = 10000
$ws.Range("N129").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O129").Value = "Región de La Araucanía"
$ws.Range("P129").Value = 833
$ws.Range("Q129").Value = 12
$ws.Range("R129").Value = "Hortaliza"
